$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held a single sample item ("Mint leaves", ...) including a
# Description value in H2. The new data set replaces that sample row and adds
# several more rows; none of the new rows use column F (HSN_Code) or column H
# (Description), so clear the stale H2 value before overwriting the rest.
$ws.Range("H2").ClearContents()

# Row 2: Tomato
$ws.Range("A2").Value = "Tomato"
$ws.Range("B2").Value = "Buy"
$ws.Range("C2").Value = "Raw Material"
$ws.Range("D2").Value = "TM_0001"
$ws.Range("E2").Value = 40.05
$ws.Range("G2").Value = 100

# Row 3: Garlic Cloves
$ws.Range("A3").Value = "Garlic Cloves"
$ws.Range("B3").Value = "Buy"
$ws.Range("C3").Value = "Raw Material"
$ws.Range("D3").Value = "GC_0002"
$ws.Range("E3").Value = 20
$ws.Range("G3").Value = 100

# Row 4: Fresh Ginger
$ws.Range("A4").Value = "Fresh Ginger"
$ws.Range("B4").Value = "Buy"
$ws.Range("C4").Value = "Raw Material"
$ws.Range("D4").Value = "FG_0003"
$ws.Range("E4").Value = 11
$ws.Range("G4").Value = 100

# Row 5: Salt
$ws.Range("A5").Value = "Salt"
$ws.Range("B5").Value = "Buy"
$ws.Range("C5").Value = "Raw Material"
$ws.Range("D5").Value = "ST_0004"
$ws.Range("E5").Value = 10
$ws.Range("G5").Value = 100

# Row 6: Kashmiri chilli powder
$ws.Range("A6").Value = "Kashmiri chilli powder"
$ws.Range("B6").Value = "Buy"
$ws.Range("C6").Value = "Raw Material"
$ws.Range("D6").Value = "KCP_0005"
$ws.Range("E6").Value = 70
$ws.Range("G6").Value = 100

# Row 7: Sunflower Oil
$ws.Range("A7").Value = "Sunflower Oil"
$ws.Range("B7").Value = "Buy"
$ws.Range("C7").Value = "Raw Material"
$ws.Range("D7").Value = "SNO_0006"
$ws.Range("E7").Value = 33
$ws.Range("G7").Value = 100

# Row 8: Melon seeds
$ws.Range("A8").Value = "Melon seeds"
$ws.Range("B8").Value = "Buy"
$ws.Range("C8").Value = "Raw Material"
$ws.Range("D8").Value = "MLS_0007"
$ws.Range("E8").Value = 32
$ws.Range("G8").Value = 100

# Row 9: Tomato Paste
$ws.Range("A9").Value = "Tomato Paste"
$ws.Range("B9").Value = "Buy"
$ws.Range("C9").Value = "Raw Material"
$ws.Range("D9").Value = "TMP_0008"
$ws.Range("E9").Value = 60
$ws.Range("G9").Value = 100

# Row 10: Sodium Tricitrate
$ws.Range("A10").Value = "Sodium Tricitrate"
$ws.Range("B10").Value = "Buy"
$ws.Range("C10").Value = "Raw Material"
$ws.Range("D10").Value = "ST_0009"
$ws.Range("E10").Value = 55
$ws.Range("G10").Value = 100

# Match the author's final selection/cursor position.
[void]$ws.Range("C11").Select()

# Page setup tweak that shipped with this commit.
$ws.PageSetup.Orientation = 1
